# Insert a new data row at row 713 (pushing existing rows 713..816 down to 714..817)
# and populate it with the new price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(713).Insert()

$ws.Cells.Item(713, 1).Value = 10
$ws.Cells.Item(713, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(713, 3).Value = "La Araucanía"
$ws.Cells.Item(713, 4).Value = 45077
$ws.Cells.Item(713, 5).Value = 9
$ws.Cells.Item(713, 6).Value = 100112032
$ws.Cells.Item(713, 7).Value = "Zapallo italiano"
$ws.Cells.Item(713, 8).Value = "Sin especificar"
$ws.Cells.Item(713, 9).Value = "Primera"
$ws.Cells.Item(713, 10).Value = 125
$ws.Cells.Item(713, 11).Value = 15000
$ws.Cells.Item(713, 12).Value = 15000
$ws.Cells.Item(713, 13).Value = 15000
$ws.Cells.Item(713, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(713, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(713, 16).Value = 300
$ws.Cells.Item(713, 17).Value = 50
$ws.Cells.Item(713, 18).Value = "Hortaliza"
